$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.670.25"
$ws.Cells.Item(2, 5).Value = "  +2.19%  "
$ws.Cells.Item(3, 4).Value = "2.553.39"
$ws.Cells.Item(3, 5).Value = "  +5.26%  "
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "570.64"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +2.52%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "151.24"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +9.06%  "
$ws.Cells.Item(7, 5).Value = "  -0.05%  "
$ws.Cells.Item(8, 5).Value = "  +0.54%  "
$ws.Cells.Item(9, 4).Value = "2.548.32"
$ws.Cells.Item(9, 5).Value = "  +5.13%  "
$ws.Cells.Item(10, 5).Value = "  +2.59%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.74"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  +0.13%  "
$ws.Cells.Item(12, 5).Value = "  +1.14%  "
$ws.Cells.Item(13, 5).Value = "  +3.09%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "28.50"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  +9.16%  "
$ws.Cells.Item(15, 4).Value = "3.010.15"
$ws.Cells.Item(15, 5).Value = "  +5.33%  "
$ws.Cells.Item(16, 4).Value = "63.567.18"
$ws.Cells.Item(16, 5).Value = "  +2.17%  "
$ws.Cells.Item(17, 5).Value = "  +2.66%  "
$ws.Cells.Item(18, 4).Value = "2.581.49"
$ws.Cells.Item(18, 5).Value = "  +5.61%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "11.71"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +4.92%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "341.36"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -1.14%  "
$ws.Cells.Item(21, 5).Value = "  +4.34%  "
$ws.Cells.Item(22, 5).Value = "  +0.74%  "
$ws.Cells.Item(23, 5).Value = "  +0.23%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "66.18"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +1.75%  "
$ws.Cells.Item(25, 5).Value = "  -0.74%  "
$ws.Cells.Item(26, 5).Value = "  +4.97%  "
$ws.Cells.Item(27, 5).Value = "  +14.19%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "8.58"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "  +5.40%  "
$ws.Cells.Item(29, 5).Value = "  +0.04%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.19"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +11.23%  "
$ws.Cells.Item(31, 4).Value = "0.0₃0832"
$ws.Cells.Item(31, 5).Value = "  +6.11%  "
$ws.Cells.Item(32, 5).Value = "  +4.82%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "177.97"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "  +3.53%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.59"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  +10.37%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "423.99"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  +12.09%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "19.22"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  +3.67%  "
$ws.Cells.Item(38, 5).Value = "  +0.91%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.79"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +6.05%  "
$ws.Cells.Item(41, 5).Value = "  +0.02%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "39.72"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +1.07%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "154.21"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  +6.63%  "
$ws.Cells.Item(44, 5).Value = "  +4.62%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "21.09"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +1.80%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.613"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  +3.99%  "
$ws.Cells.Item(47, 5).Value = "  +2.58%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0971"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = "  +1.97%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0239"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +7.65%  "
$ws.Cells.Item(50, 5).Value = "  +4.92%  "
$ws.Cells.Item(51, 5).Value = "  +8.18%  "
